# configs/user_position.xlsx — "bs env for sicmdp"
# Shrink the sample data from 7 data rows down to a single data row (row 2),
# with new x/y values, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets new values (A2: 2000 -> 200, B2: 3500 -> 130). C2 (1.5) is left as-is.
$ws.Range("A2").Value = 200
$ws.Range("B2").Value = 130

# Rows 3-8 are removed entirely, shrinking the sheet's used range to A1:C2.
$ws.Range("A3:C8").EntireRow.Delete()

# The active selection moves to E5.
$ws.Range("E5").Select()
